# Some fix Print list tabs
# Add a new worksheet "New Лист" after "Лист1" and populate it with a
# small print list (One, Two, Three, Four, Five) in column A.

$wb = $excel.ActiveWorkbook

$firstSheet = $wb.Worksheets.Item(1)

# Insert the new worksheet right after the first ("Лист1") sheet.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $firstSheet)
$ws.Name = "New Лист"

$ws.Range("A1").Value = "One"
$ws.Range("A2").Value = "Two"
$ws.Range("A3").Value = "Three"
$ws.Range("A4").Value = "Four"
$ws.Range("A5").Value = "Five"
